$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the new weekly record (shifts existing rows 507:539 down to 509:541)
$ws.Rows("507:508").Insert()

# Row 507 - "Primera" quality record for the new week (2022-01-24)
$ws.Range("A507").Value = 3
$ws.Range("B507").Value = "Femacal de La Calera"
$ws.Range("C507").Value = "Coquimbo"
$ws.Range("D507").Value = 44585
$ws.Range("E507").Value = 5
$ws.Range("F507").Value = 100114014
$ws.Range("G507").Value = "Betarraga"
$ws.Range("H507").Value = "Sin especificar"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 2700
$ws.Range("K507").Value = 550
$ws.Range("L507").Value = 600
$ws.Range("M507").Value = 576
$ws.Range("N507").Value = "$/paquete 4 unidades"
$ws.Range("O507").Value = "Provincia de Quillota"
$ws.Range("P507").Value = 144
$ws.Range("Q507").Value = 4
$ws.Range("R507").Value = "Hortaliza"

# Row 508 - "Segunda" quality record for the new week (2022-01-24)
$ws.Range("A508").Value = 3
$ws.Range("B508").Value = "Femacal de La Calera"
$ws.Range("C508").Value = "Coquimbo"
$ws.Range("D508").Value = 44585
$ws.Range("E508").Value = 5
$ws.Range("F508").Value = 100114014
$ws.Range("G508").Value = "Betarraga"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Segunda"
$ws.Range("J508").Value = 1500
$ws.Range("K508").Value = 400
$ws.Range("L508").Value = 400
$ws.Range("M508").Value = 400
$ws.Range("N508").Value = "$/paquete 4 unidades"
$ws.Range("O508").Value = "Provincia de Quillota"
$ws.Range("P508").Value = 100
$ws.Range("Q508").Value = 4
$ws.Range("R508").Value = "Hortaliza"
